$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4-12, keeping only the header row and the first two data rows
$ws.Range("A4:C12").EntireRow.Delete() | Out-Null

# Update the remaining data rows with the new book/author values
$ws.Range("B2").Value = "A Love Supreme"
$ws.Range("C2").Value = "John Coltrane"
$ws.Range("B3").Value = "VALIS"
$ws.Range("C3").Value = "Philip K. Dick"
